# Update to plasticity AIC table
$wb = $excel.ActiveWorkbook

# --- Sheet "Table S1 - Plasticity AIC" ---
$ws1 = $wb.Worksheets.Item("Table S1 - Plasticity AIC")

# Delete the erroneous duplicate "pstr" row (old row 9, model formula had a stray
# "reef environment * pCO2 * temperature + (1 | colony)7" entry). This removes the
# bad model row and shifts all following rows up by one.
$ws1.Rows.Item(9).Delete()

# Update rounded Marginal/Conditional R2 values (columns F and G) for all remaining
# data rows (rows 2-17 after the deletion above).
$ws1.Range("F2").Value = 0.506
$ws1.Range("G2").Value = 0.322

$ws1.Range("F3").Value = 0.545
$ws1.Range("G3").Value = 0.366

$ws1.Range("F4").Value = 0.512
$ws1.Range("G4").Value = 0.329

$ws1.Range("F5").Value = 0.442
$ws1.Range("G5").Value = 0.254

$ws1.Range("F6").Value = 0.37
$ws1.Range("G6").Value = 0.088

$ws1.Range("F7").Value = 0.442
$ws1.Range("G7").Value = 0.253

$ws1.Range("F8").Value = 0.397

$ws1.Range("F9").Value = 0.309
$ws1.Range("G9").Value = 0.261

$ws1.Range("F10").Value = 0.278
$ws1.Range("G10").Value = 0.238

$ws1.Range("F11").Value = 0.232
$ws1.Range("G11").Value = 0.188

$ws1.Range("F12").Value = 0.521

$ws1.Range("F13").Value = 0.522

$ws1.Range("F14").Value = 0.527
$ws1.Range("G14").Value = 0.199

$ws1.Range("F15").Value = 0.499
$ws1.Range("G15").Value = 0.174

$ws1.Range("F16").Value = 0.485
$ws1.Range("G16").Value = 0.147

$ws1.Range("F17").Value = 0.5
$ws1.Range("G17").Value = 0.174

# --- Sheet "Table S5 - HostVsymb PERMANOVA" ---
$ws5 = $wb.Worksheets.Item("Table S5 - HostVsymb PERMANOVA")
$ws5.Range("F2").Value = 0.74217
$ws5.Range("F4").Value = 0.56762
